# miata_dash_v1.0 - "made changes to the excel formulas"
#
# Adds three volt/temp/fuel/oil-psi conversion tables (rows 19-42) below the
# existing pin-out table, with their own header banners, alternating
# red/blue/cyan/olive highlight styles, and live formulas that convert a
# sender's measured voltage into a calibrated gauge reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths (A:E) - leave F.. untouched
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 14.5066666667
$ws.Columns("B").ColumnWidth = 11.7266666667
$ws.Columns("C").ColumnWidth = 12.9766666667
$ws.Columns("D").ColumnWidth = 10.6866666667
$ws.Columns("E").ColumnWidth = 15.8966666667

# ---------------------------------------------------------------------------
# Style helpers -- reproduce the fill / font / numberformat combinations the
# author applied through Format Cells.  Interior.Color supplies the fgColor,
# Interior.PatternColor the (cosmetic) bgColor; both take BGR-packed ints.
# ---------------------------------------------------------------------------
function Style-Range($rng, $fg, $bg, $bold, $white, $numfmt, $center) {
    if ($fg -ne $null) {
        $rng.Interior.Color = $fg
        $rng.Interior.PatternColor = $bg
    }
    $rng.Font.Bold = $bold
    if ($white) {
        $rng.Font.Color = 16777215
    } else {
        $rng.Font.Color = 0
    }
    $rng.NumberFormat = $numfmt
    if ($center) {
        $rng.HorizontalAlignment = -4108
    }
}

function Style-Title($rng)      { Style-Range $rng 2013534  32896   $false $false "General" $true }
function Style-TitleE($rng)     { Style-Range $rng 2013534  32896   $false $false "#.00"    $true }
function Style-PlainC($rng)     { Style-Range $rng $null    $null   $false $false "General" $true }
function Style-Num00C($rng)     { Style-Range $rng $null    $null   $false $false "0.00"    $true }
function Style-NumPtC($rng)     { Style-Range $rng $null    $null   $false $false "#.00"    $true }
function Style-BoldRedGen($rng) { Style-Range $rng 66047    255     $true  $true  "General" $true }
function Style-BoldBluePt($rng) { Style-Range $rng 11567961 8421504 $true  $true  "#.00"    $true }
function Style-Bold00Red($rng)  { Style-Range $rng 66047    255     $true  $true  "0.00"    $true }
function Style-CyanGen($rng)    { Style-Range $rng 16561416 13421619 $false $false "General" $true }
function Style-Olive00($rng)    { Style-Range $rng 2013534  32896   $false $false "0.00"    $true }
function Style-WhiteRedGen($rng){ Style-Range $rng 66047    255     $false $true  "General" $true }
function Style-Plain00Gen($rng) { Style-Range $rng $null    $null   $false $false "0.00"    $false }

# ---------------------------------------------------------------------------
# Row 19 - "Temp " banner (re-uses the existing "Temp " label text)
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "Temp "
Style-Title  $ws.Range("A19:D19")
Style-TitleE $ws.Range("E19")

# Row 20 - headers
$ws.Range("A20").Value = "NORM OP TEMP"
$ws.Range("B20").Value = "VOLTAGE"
$ws.Range("C20").Value = "START VOLTS"
$ws.Range("D20").Value = "VOLT DIFF"
$ws.Range("E20").Value = "CONVERT RATIO"
Style-PlainC $ws.Range("A20")
Style-PlainC $ws.Range("B20")
Style-Num00C $ws.Range("C20")
Style-PlainC $ws.Range("D20")
Style-NumPtC $ws.Range("E20")

# Row 21 - calibration data + formulas
$ws.Range("A21").Value = 195
$ws.Range("B21").Value = 5.13
$ws.Range("C21").Value = 10
$ws.Range("D21").Formula = "=C21-B21"
$ws.Range("E21").Formula = "=(A21*0.01) / D21"
Style-PlainC     $ws.Range("A21")
Style-PlainC     $ws.Range("B21")
Style-Num00C     $ws.Range("C21")
Style-BoldRedGen $ws.Range("D21")
Style-BoldBluePt $ws.Range("E21")

# Rows 22-23 - blank spacer rows (style only)
foreach ($r in 22..23) {
    Style-PlainC $ws.Range("A$r")
    Style-PlainC $ws.Range("B$r")
    Style-Num00C $ws.Range("C$r")
    Style-PlainC $ws.Range("D$r")
    Style-NumPtC $ws.Range("E$r")
}

# Row 24 - 2nd headers
$ws.Range("A24").Value = "CURRENT VOLT"
$ws.Range("B24").Value = "START VOLT"
$ws.Range("C24").Value = "VOLT DIFF"
$ws.Range("D24").Value = "TEMP"
$ws.Range("E24").Value = "       "
Style-PlainC $ws.Range("A24")
Style-PlainC $ws.Range("B24")
Style-Num00C $ws.Range("C24")
Style-PlainC $ws.Range("D24")
Style-NumPtC $ws.Range("E24")

# Row 25 - live readout
$ws.Range("A25").Value = 5.13
$ws.Range("B25").Formula = "=C21"
$ws.Range("C25").Formula = "=B25-A25"
$ws.Range("D25").Formula = "=(C25*100) * E21"
Style-PlainC  $ws.Range("A25")
Style-Num00C  $ws.Range("B25")
Style-Bold00Red $ws.Range("C25")
Style-CyanGen $ws.Range("D25")
Style-NumPtC  $ws.Range("E25")

# Row 26 - blank spacer
Style-PlainC $ws.Range("A26")
Style-PlainC $ws.Range("B26")
Style-Num00C $ws.Range("C26")
Style-PlainC $ws.Range("D26")
Style-NumPtC $ws.Range("E26")

# ---------------------------------------------------------------------------
# Row 27 - "FUEL" banner
# ---------------------------------------------------------------------------
$ws.Range("A27").Value = "FUEL"
Style-Title  $ws.Range("A27:B27")
Style-Olive00 $ws.Range("C27")
Style-Title  $ws.Range("D27")
Style-TitleE $ws.Range("E27")

# Row 28 - headers
$ws.Range("A28").Value = "1 GALLON"
$ws.Range("B28").Value = "VOLTAGE"
$ws.Range("C28").Value = "START VOLTS"
$ws.Range("D28").Value = "VOLT DIFF"
$ws.Range("E28").Value = "CONVERT RATIO"
Style-PlainC $ws.Range("A28")
Style-PlainC $ws.Range("B28")
Style-Num00C $ws.Range("C28")
Style-PlainC $ws.Range("D28")
Style-NumPtC $ws.Range("E28")

# Row 29 - calibration data + formulas
$ws.Range("A29").Value = 1
$ws.Range("B29").Value = 0.83333
$ws.Range("C29").Value = 10
$ws.Range("D29").Formula = "=C29-B29"
$ws.Range("E29").Formula = "=(A29*0.01) / D29"
Style-PlainC       $ws.Range("A29")
Style-PlainC       $ws.Range("B29")
Style-Num00C       $ws.Range("C29")
Style-WhiteRedGen  $ws.Range("D29")
Style-BoldBluePt   $ws.Range("E29")

# Rows 30-31 - blank spacer rows
foreach ($r in 30..31) {
    Style-PlainC $ws.Range("A$r")
    Style-PlainC $ws.Range("B$r")
    Style-Num00C $ws.Range("C$r")
    Style-PlainC $ws.Range("D$r")
    Style-NumPtC $ws.Range("E$r")
}

# Row 32 - 2nd headers
$ws.Range("A32").Value = "CURRENT VOLT"
$ws.Range("B32").Value = "START VOLT"
$ws.Range("C32").Value = "VOLT DIFF"
$ws.Range("D32").Value = "FUEL"
$ws.Range("E32").Value = "       "
Style-PlainC $ws.Range("A32")
Style-PlainC $ws.Range("B32")
Style-Num00C $ws.Range("C32")
Style-PlainC $ws.Range("D32")
Style-NumPtC $ws.Range("E32")

# Row 33 - live readout
$ws.Range("A33").Value = 5
$ws.Range("B33").Formula = "=C29"
$ws.Range("C33").Formula = "=B33-A33"
$ws.Range("D33").Formula = "=(C33*100) * E29"
Style-PlainC    $ws.Range("A33")
Style-Num00C    $ws.Range("B33")
Style-Bold00Red $ws.Range("C33")
Style-CyanGen   $ws.Range("D33")
Style-NumPtC    $ws.Range("E33")

# Row 34 - blank spacer
Style-PlainC $ws.Range("A34")
Style-PlainC $ws.Range("B34")
Style-Num00C $ws.Range("C34")
Style-PlainC $ws.Range("D34")
Style-NumPtC $ws.Range("E34")

# ---------------------------------------------------------------------------
# Row 35 - "OIL PSI" banner
# ---------------------------------------------------------------------------
$ws.Range("A35").Value = "OIL PSI"
Style-Title  $ws.Range("A35:B35")
Style-Olive00 $ws.Range("C35")
Style-Title  $ws.Range("D35")
Style-TitleE $ws.Range("E35")

# Row 36 - headers
$ws.Range("A36").Value = "30 PSI"
$ws.Range("B36").Value = "VOLTAGE"
$ws.Range("C36").Value = "START VOLTS"
$ws.Range("D36").Value = "VOLT DIFF"
$ws.Range("E36").Value = "CONVERT RATIO"
Style-PlainC $ws.Range("A36")
Style-PlainC $ws.Range("B36")
Style-Num00C $ws.Range("C36")
Style-PlainC $ws.Range("D36")
Style-NumPtC $ws.Range("E36")

# Row 37 - calibration data + formulas
$ws.Range("A37").Value = 30
$ws.Range("B37").Value = 5.13
$ws.Range("C37").Value = 10
$ws.Range("D37").Formula = "=C37-B37"
$ws.Range("E37").Formula = "=(A37*0.01) / D37"
Style-PlainC      $ws.Range("A37")
Style-PlainC      $ws.Range("B37")
Style-Num00C      $ws.Range("C37")
Style-WhiteRedGen $ws.Range("D37")
Style-BoldBluePt  $ws.Range("E37")

# Rows 38-39 - blank spacer rows
foreach ($r in 38..39) {
    Style-PlainC $ws.Range("A$r")
    Style-PlainC $ws.Range("B$r")
    Style-Num00C $ws.Range("C$r")
    Style-PlainC $ws.Range("D$r")
    Style-NumPtC $ws.Range("E$r")
}

# Row 40 - 2nd headers
$ws.Range("A40").Value = "CURRENT VOLT"
$ws.Range("B40").Value = "START VOLT"
$ws.Range("C40").Value = "VOLT DIFF"
$ws.Range("D40").Value = "PSI"
$ws.Range("E40").Value = "       "
Style-PlainC $ws.Range("A40")
Style-PlainC $ws.Range("B40")
Style-Num00C $ws.Range("C40")
Style-PlainC $ws.Range("D40")
Style-NumPtC $ws.Range("E40")

# Row 41 - live readout
$ws.Range("A41").Value = 8.13
$ws.Range("B41").Formula = "=C37"
$ws.Range("C41").Formula = "=B41-A41"
$ws.Range("D41").Formula = "=(C41*100) * E37"
Style-PlainC    $ws.Range("A41")
Style-Num00C    $ws.Range("B41")
Style-Bold00Red $ws.Range("C41")
Style-CyanGen   $ws.Range("D41")
Style-NumPtC    $ws.Range("E41")

# Row 42 - trailing spacer cell
Style-Plain00Gen $ws.Range("C42")

# ---------------------------------------------------------------------------
# Leave the selection where the author left it
# ---------------------------------------------------------------------------
$ws.Range("H23").Select()
